$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Add a new labelled table of operator pairs in rows 1-9 (mirrors the
# existing table at rows 14-22, but uses plain numeric indices 0..63
# instead of the symbolic x1..x8 markers), plus a label in O2 pointing
# back at it.
# -----------------------------------------------------------------------

$opNames = @("streamFilter","streamMap","streamFilterAcc","streamScan","streamWindow","streamExpand","streamJoin","streamMerge")

# Row 1: header row (same formatting as row 14's header) -----------------
$ws.Range("E14:M14").Copy()
$ws.Range("E1:M1").PasteSpecial(-4122)
$ws.Rows(1).RowHeight = 74.35

for ($c = 6; $c -le 13; $c++) {
  $ws.Cells.Item(1, $c).Value = $opNames[$c - 6]
}

# Rows 2-9: operator-name column + numbered grid --------------------------
$ws.Range("E15").Copy()
$ws.Range("E2:E9").PasteSpecial(-4122)
$ws.Rows("2:9").RowHeight = 12.8

for ($r = 2; $r -le 9; $r++) {
  $ws.Cells.Item($r, 5).Value = $opNames[$r - 2]
  for ($c = 6; $c -le 13; $c++) {
    $ws.Cells.Item($r, $c).Value = ($r - 2) * 8 + ($c - 6)
  }
}

# Label pointing at the new table ------------------------------------------
$ws.Range("O2").Value = [char]0x2190 + "labelling unique operator pairs"

# -----------------------------------------------------------------------
# View/selection bookkeeping: mirror the authored change to sheetView,
# which re-selects O3 as the active cell (and scrolls back to the top).
# -----------------------------------------------------------------------
$ws.Range("O3").Select() | Out-Null
